# MITHEx_inputs.xlsx update:
# - Update iteration method and improve cycle efficiency calculations.
# - Add file for iterating on parameters.

$wb = $excel.ActiveWorkbook

$wsPlant  = $wb.Worksheets.Item("Plant Description")
$wsHX     = $wb.Worksheets.Item("HX Parameters")
$wsCycle  = $wb.Worksheets.Item("Cycle Parameters")
$wsInput  = $wb.Worksheets.Item("Input options")

# ---- Plant Description sheet ----
# Primary fluid -> Air instead of CarbonDioxide
$wsPlant.Range("B7").Value = "Air"
# Secondary Cold Temperature (C)
$wsPlant.Range("B9").Value = 250
# Secondary Pressure (kPa)
$wsPlant.Range("B11").Value = 6000
$wsPlant.Range("B6").Select()

# ---- HX Parameters sheet ----
# Plate thickness (m)
$wsHX.Range("B2").Value = 0.002
# HX length lower bound (m)
$wsHX.Range("B7").Value = 1
# HX length upper bound (m)
$wsHX.Range("B8").Value = 10
$wsHX.Range("B9").Select()

# ---- Cycle Parameters sheet ----
# Pump/Compressor Efficiency
$wsCycle.Range("B1").Value = 0.88
$wsCycle.Range("C1").Value = "Efficiency of the pump/compressor in the power concersion cycle"
# Turbine Efficiency
$wsCycle.Range("B2").Value = 0.88
$wsCycle.Range("C2").Value = "Efficiency of the turbine"
# Compression Ratio
$wsCycle.Range("C3").Value = "Ratio between high and low pressures"

# New row 4: Secondary Minimum Temperature (C) - copy formatting from row 3
# so the new cells pick up the same borders/fill as the rest of the table.
$wsCycle.Range("A3").Copy()
$wsCycle.Range("A4").PasteSpecial(-4122)
$wsCycle.Range("A4").Value = "Secondary Minimum Temperature (C)"

$wsCycle.Range("B3").Copy()
$wsCycle.Range("B4").PasteSpecial(-4122)

$wsCycle.Range("C3").Copy()
$wsCycle.Range("C4").PasteSpecial(-4122)
$wsCycle.Range("C4").Value = "(Optional) Allows for simplified modeling of a recuperator where the pump inlet temperature is specified and used for pump calculations instead of the secondary cold temperature"
$wsCycle.Range("C4").WrapText = $true

$wsCycle.Rows.Item(4).RowHeight = 64

$wsCycle.Columns.Item(1).ColumnWidth = 29.833333333333332
$wsCycle.Columns.Item(3).ColumnWidth = 38.333333333333336

$wsCycle.Range("B4").Select()

# Make HX Parameters the active sheet/tab
$wsHX.Activate()
$wsHX.Range("B9").Select()
